$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Column A (combination) text for rows 4-19 (row 3 "PTK2" unchanged) ---
# --- and Column I (Group) text for rows 3-19: "Pass 2" -> "17-gene set" ---
$ws.Range("A3").Value2 = 'PTK2'
$ws.Range("A4").Value2 = 'PTK2, TNFRSF13B'
$ws.Range("A5").Value2 = 'PTK2, TNFRSF13B, UBE2G1'
$ws.Range("A6").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3'
$ws.Range("A7").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2'
$ws.Range("A8").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2'
$ws.Range("A9").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A'
$ws.Range("A10").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3'
$ws.Range("A11").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3, GTSE1'
$ws.Range("A12").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3, GTSE1, TRIQK'
$ws.Range("A13").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3, GTSE1, TRIQK, ADPGK'
$ws.Range("A14").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3, GTSE1, TRIQK, ADPGK, GDPGP1'
$ws.Range("A15").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3, GTSE1, TRIQK, ADPGK, GDPGP1, CCR5'
$ws.Range("A16").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3, GTSE1, TRIQK, ADPGK, GDPGP1, CCR5, ZNF628'
$ws.Range("A17").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3, GTSE1, TRIQK, ADPGK, GDPGP1, CCR5, ZNF628, NABP1'
$ws.Range("A18").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3, GTSE1, TRIQK, ADPGK, GDPGP1, CCR5, ZNF628, NABP1, GEMIN5'
$ws.Range("A19").Value2 = 'PTK2, TNFRSF13B, UBE2G1, SENP3, UQCC2, TAP2, MON1A, NFE2L3, GTSE1, TRIQK, ADPGK, GDPGP1, CCR5, ZNF628, NABP1, GEMIN5, GDAP2'

for ($r = 3; $r -le 19; $r++) {
    $ws.Range("I$r").Value2 = '17-gene set'
}

# --- Update recomputed numeric cells (Fold1..Fold5, Average) ---
$ws.Range("C8").Value2 = 0.93188854489164097
$ws.Range("G8").Value2 = 0.950236863216947
$ws.Range("C10").Value2 = 0.95356037151702799
$ws.Range("G10").Value2 = 0.96041412820188798
$ws.Range("B11").Value2 = 0.97278911564625803
$ws.Range("C11").Value2 = 0.97523219814241502
$ws.Range("D11").Value2 = 0.94140625
$ws.Range("E11").Value2 = 0.97083333333333299
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 0.97205217942440103
$ws.Range("B12").Value2 = 0.98979591836734704
$ws.Range("C12").Value2 = 0.95665634674922595
$ws.Range("E12").Value2 = 0.98750000000000004
$ws.Range("G12").Value2 = 0.97975920302331498
$ws.Range("C13").Value2 = 0.95975232198142402
$ws.Range("D13").Value2 = 0.9609375
$ws.Range("E13").Value2 = 0.99583333333333302
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 0.98330463106295196
$ws.Range("C14").Value2 = 0.98142414860681104
$ws.Range("D14").Value2 = 0.97265625
$ws.Range("E14").Value2 = 0.97083333333333299
$ws.Range("G14").Value2 = 0.98362220217034202
$ws.Range("B15").Value2 = 0.98979591836734704
$ws.Range("C15").Value2 = 0.99071207430340602
$ws.Range("D15").Value2 = 0.9765625
$ws.Range("E15").Value2 = 0.97916666666666696
$ws.Range("G15").Value2 = 0.98724743186748398
$ws.Range("B16").Value2 = 0.99659863945578198
$ws.Range("C16").Value2 = 0.99690402476780204
$ws.Range("D16").Value2 = 0.97265625
$ws.Range("G16").Value2 = 0.99323178284471703
$ws.Range("B17").Value2 = 1
$ws.Range("C17").Value2 = 1
$ws.Range("E17").Value2 = 0.99166666666666703
$ws.Range("G17").Value2 = 0.99520833333333303
$ws.Range("B18").Value2 = 1
$ws.Range("D18").Value2 = 0.9921875
$ws.Range("G18").Value2 = 0.99843749999999998
$ws.Range("B19").Value2 = 0.99319727891156495
$ws.Range("D19").Value2 = 1
$ws.Range("G19").Value2 = 0.99863945578231295

# --- Row heights: rows 2-19 grow from 20.1 to 24.95 to match the header row ---
$ws.Range("A2:A19").EntireRow.RowHeight = 24.95

# --- New bestFit width on column I (Group) after the longer "17-gene set" label ---
$ws.Columns(9).AutoFit()
